# Update countries & provincias Spain
# - Israel's case count overtook Noruega and Brasil, so it moves up one rank
#   (row 20), pushing Noruega to row 21 and Brasil to row 22.
# - Jordania's case count overtook Kuwait, so it moves up one rank (row 85),
#   pushing Kuwait to row 86.
# - Refresh the "Casos totales / Nuevos casos / Casos activos / Recuperados /
#   Casos criticos / Muertes hoy / Muertes" figures for several countries.
# - Bump the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp row (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Marzo de 2020 a las 19:50"

# --- Estados Unidos (row 4): values refreshed, ranking unchanged ---
$ws.Cells.Item(4, 2).Value = 152631
$ws.Cells.Item(4, 3).Value = 9140
$ws.Cells.Item(4, 4).Value = 5211
$ws.Cells.Item(4, 5).Value = 144603
$ws.Cells.Item(4, 6).Value = 3402
$ws.Cells.Item(4, 7).Value = 234
$ws.Cells.Item(4, 8).Value = 2817

# --- Francia (row 10): values refreshed, ranking unchanged ---
$ws.Cells.Item(10, 4).Value = 7927
$ws.Cells.Item(10, 5).Value = 29223
$ws.Cells.Item(10, 6).Value = 5056
$ws.Cells.Item(10, 7).Value = 418
$ws.Cells.Item(10, 8).Value = 3024

# --- Rows 20-22: Israel overtakes Noruega and Brasil ---
# Row 20 becomes Israel, with its refreshed totals.
$ws.Cells.Item(20, 1).Value = "Israel"
$ws.Cells.Item(20, 2).Value = 4695
$ws.Cells.Item(20, 3).Value = 448
$ws.Cells.Item(20, 4).Value = 134
$ws.Cells.Item(20, 5).Value = 4545
$ws.Cells.Item(20, 6).Value = 66
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 16

# Row 21 becomes Noruega (its figures are unchanged, just shifted down a row).
$ws.Cells.Item(21, 1).Value = "Noruega"
$ws.Cells.Item(21, 2).Value = 4445
$ws.Cells.Item(21, 3).Value = 161
$ws.Cells.Item(21, 4).Value = 12
$ws.Cells.Item(21, 5).Value = 4401
$ws.Cells.Item(21, 6).Value = 97
$ws.Cells.Item(21, 7).Value = 6
$ws.Cells.Item(21, 8).Value = 32

# Row 22 becomes Brasil (its figures are unchanged, just shifted down a row).
$ws.Cells.Item(22, 1).Value = "Brasil"
$ws.Cells.Item(22, 2).Value = 4371
$ws.Cells.Item(22, 3).Value = 115
$ws.Cells.Item(22, 4).Value = 120
$ws.Cells.Item(22, 5).Value = 4110
$ws.Cells.Item(22, 6).Value = 296
$ws.Cells.Item(22, 7).Value = 5
$ws.Cells.Item(22, 8).Value = 141

# --- Sudafrica (row 42): values refreshed, ranking unchanged ---
$ws.Cells.Item(42, 2).Value = 1326
$ws.Cells.Item(42, 3).Value = 46
$ws.Cells.Item(42, 5).Value = 1292
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = 3

# --- Rows 85-86: Jordania overtakes Kuwait ---
# Row 85 becomes Jordania, with its refreshed totals.
$ws.Cells.Item(85, 1).Value = "Jordania"
$ws.Cells.Item(85, 2).Value = 268
$ws.Cells.Item(85, 3).Value = 9
$ws.Cells.Item(85, 4).Value = 26
$ws.Cells.Item(85, 5).Value = 237
$ws.Cells.Item(85, 6).Value = 3
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 5

# Row 86 becomes Kuwait (its figures are unchanged, just shifted down a row).
$ws.Cells.Item(86, 1).Value = "Kuwait"
$ws.Cells.Item(86, 2).Value = 266
$ws.Cells.Item(86, 3).Value = 11
$ws.Cells.Item(86, 4).Value = 72
$ws.Cells.Item(86, 5).Value = 194
$ws.Cells.Item(86, 6).Value = 13
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
